$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("placesToGo")

# Update the "Miami" entry to "Kansas City"
$ws.Range("A4").Value = "Kansas City"

# Move selection to A5, mirroring the user pressing Enter after editing A4
$ws.Range("A5").Select()
